# "test partage xls 02"
#
# Alex opens the shared workbook, types a second value ("test2") into C1,
# and saves a personal custom view of the sheet before leaving the
# selection on C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value in C1 - grows the shared-string table ("test2"), the used
# range (A1 -> A1:C1) and the row span (1:1 -> 1:3).
$ws.Range("C1").Value = "test2"

# Park the cursor on A2 (the cell the custom view will remember), then
# save Alex's personal custom view of this workbook/sheet.
$ws.Range("A2").Select() | Out-Null
$wb.CustomViews.Add("Alex - Affichage personnalisé") | Out-Null

# Finally Alex leaves the selection on C2, next to the value just typed.
$ws.Range("C2").Select() | Out-Null
